$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.353.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.203.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.60%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.83%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.49%  "
$ws.Range("E14").Value = "  -4.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.536.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.854"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -13.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -11.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.204.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.285.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -11.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0934"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.79%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.86%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.65%  "
$ws.Range("E26").Value = "  -9.79%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.19%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -15.28%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0861"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -13.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.06%  "
$ws.Range("E37").Value = "  -7.71%  "
$ws.Range("E38").Value = "  -8.95%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.53%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -13.45%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0315"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.03%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.762.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "87.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.202"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -14.58%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "75.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -12.17%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "59.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -15.90%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.44%  "
